$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.049841917968287
$ws.Range("D2").Value = 1.053400617270445
$ws.Range("E2").Value = 1.046930530253106
$ws.Range("F2").Value = 1.062666392331433
$ws.Range("I2").Value = 1.0435727354514
$ws.Range("J2").Value = 1.054878226620457
$ws.Range("K2").Value = 1.056146689430899
$ws.Range("L2").Value = 1.049694587561333
$ws.Range("M2").Value = 1.06538712101325

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.051358347618561
$ws.Range("D3").Value = 1.054557327506738
$ws.Range("E3").Value = 1.048242175690228
$ws.Range("F3").Value = 1.063971957731935
$ws.Range("I3").Value = 1.043978837764573
$ws.Range("J3").Value = 1.056041327173583
$ws.Range("K3").Value = 1.057115556701128
$ws.Range("L3").Value = 1.050816677690361
$ws.Range("M3").Value = 1.066506313909523

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.052338173658039
$ws.Range("D4").Value = 1.055304362615863
$ws.Range("E4").Value = 1.049089793449873
$ws.Range("F4").Value = 1.06481555939605
$ws.Range("I4").Value = 1.044239353842964
$ws.Range("J4").Value = 1.056792135404184
$ws.Range("K4").Value = 1.057740469206848
$ws.Range("L4").Value = 1.05154109995509
$ws.Range("M4").Value = 1.067228772077114

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.052749763570394
$ws.Range("D5").Value = 1.055618077640462
$ws.Range("E5").Value = 1.049445873179391
$ws.Range("F5").Value = 1.065169930766175
$ws.Range("I5").Value = 1.04434833609332
$ws.Range("J5").Value = 1.057107350766832
$ws.Range("K5").Value = 1.058002705792261
$ws.Range("L5").Value = 1.051845258430361
$ws.Range("M5").Value = 1.067532083104962

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.052818852254557
$ws.Range("D6").Value = 1.05567073205791
$ws.Range("E6").Value = 1.049505645521732
$ws.Range("F6").Value = 1.065229415083847
$ws.Range("I6").Value = 1.044366603165147
$ws.Range("J6").Value = 1.057160252083045
$ws.Range("K6").Value = 1.058046708646905
$ws.Range("L6").Value = 1.051896305336695
$ws.Range("M6").Value = 1.067582986436317

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.052343674625897
$ws.Range("D7").Value = 1.055308555815779
$ws.Range("E7").Value = 1.049094552414469
$ws.Range("F7").Value = 1.064820295612509
$ws.Range("I7").Value = 1.044240812183228
$ws.Range("J7").Value = 1.056796348988228
$ws.Range("K7").Value = 1.057743975092326
$ws.Range("L7").Value = 1.05154516565622
$ws.Range("M7").Value = 1.067232826540468

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.050354697972817
$ws.Range("D8").Value = 1.053791832624423
$ws.Range("E8").Value = 1.047374039085284
$ws.Range("F8").Value = 1.06310786311509
$ws.Range("I8").Value = 1.043710448735698
$ws.Range("J8").Value = 1.055271676980944
$ws.Range("K8").Value = 1.056474541392767
$ws.Range("L8").Value = 1.050074146931831
$ws.Range("M8").Value = 1.065765719556167

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.046838767276416
$ws.Range("D9").Value = 1.051107979521698
$ws.Range("E9").Value = 1.04433356105459
$ws.Range("F9").Value = 1.060081014239943
$ws.Range("I9").Value = 1.042758479484685
$ws.Range("J9").Value = 1.052571012620833
$ws.Range("K9").Value = 1.054222048159246
$ws.Range("L9").Value = 1.047469191202956
$ws.Range("M9").Value = 1.06316698715619

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.044486878235401
$ws.Range("D10").Value = 1.049310919235608
$ws.Range("E10").Value = 1.042300359444692
$ws.Range("F10").Value = 1.058056511637484
$ws.Range("I10").Value = 1.042111995430485
$ws.Range("J10").Value = 1.050760813347251
$ws.Range("K10").Value = 1.052709623999339
$ws.Range("L10").Value = 1.045723597580538
$ws.Range("M10").Value = 1.061425112625559

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.043466496473812
$ws.Range("D11").Value = 1.048530848934371
$ws.Range("E11").Value = 1.041418406056197
$ws.Range("F11").Value = 1.057178238213671
$ws.Range("I11").Value = 1.041829221596281
$ws.Range("J11").Value = 1.049974585501838
$ws.Range("K11").Value = 1.052052111809339
$ws.Range("L11").Value = 1.044965538224853
$ws.Range("M11").Value = 1.060668566129912

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.043087171848374
$ws.Range("D12").Value = 1.048240800025716
$ws.Range("E12").Value = 1.041090567379339
$ws.Range("F12").Value = 1.05685175382197
$ws.Range("I12").Value = 1.041723757055159
$ws.Range("J12").Value = 1.049682178191638
$ws.Range("K12").Value = 1.051807482955586
$ws.Range("L12").Value = 1.044683623398943
$ws.Range("M12").Value = 1.060387198959075

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.043168552417708
$ws.Range("D13").Value = 1.048303029994167
$ws.Range("E13").Value = 1.04116090094864
$ws.Range("F13").Value = 1.056921797456787
$ws.Range("I13").Value = 1.041746399049078
$ws.Range("J13").Value = 1.049744917281751
$ws.Range("K13").Value = 1.051859974845115
$ws.Range("L13").Value = 1.044744110477337
$ws.Range("M13").Value = 1.060447569204346

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.043435147747131
$ws.Range("D14").Value = 1.048506879475096
$ws.Range("E14").Value = 1.041391311784635
$ws.Range("F14").Value = 1.057151256142727
$ws.Range("I14").Value = 1.041820512647293
$ws.Range("J14").Value = 1.049950422554128
$ws.Range("K14").Value = 1.052031898894916
$ws.Range("L14").Value = 1.044942241986386
$ws.Range("M14").Value = 1.060645315450163

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.043599364652964
$ws.Range("D15").Value = 1.048632438459045
$ws.Range("E15").Value = 1.041533243238368
$ws.Range("F15").Value = 1.057292599285905
$ws.Range("I15").Value = 1.041866119451921
$ws.Range("J15").Value = 1.050076992263248
$ws.Range("K15").Value = 1.052137773865267
$ws.Range("L15").Value = 1.045064272389438
$ws.Range("M15").Value = 1.060767106597331

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.044554554728185
$ws.Range("D16").Value = 1.049362648738496
$ws.Range("E16").Value = 1.042358858245775
$ws.Range("F16").Value = 1.058114764419106
$ws.Range("I16").Value = 1.042130702085243
$ws.Range("J16").Value = 1.050812941566493
$ws.Range("K16").Value = 1.052753205191828
$ws.Range("L16").Value = 1.045773860461812
$ws.Range("M16").Value = 1.061475273079571

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.045153178763843
$ws.Range("D17").Value = 1.049820169127762
$ws.Range("E17").Value = 1.042876321165435
$ws.Range("F17").Value = 1.058630040390084
$ws.Range("I17").Value = 1.042295904953317
$ws.Range("J17").Value = 1.051273935915162
$ws.Range("K17").Value = 1.053138542974237
$ws.Range("L17").Value = 1.046218371231635
$ws.Range("M17").Value = 1.061918866269384

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.045502153854807
$ws.Range("D18").Value = 1.050086847086184
$ws.Range("E18").Value = 1.043177998135681
$ws.Range("F18").Value = 1.058930433350786
$ws.Range("I18").Value = 1.042391990971234
$ws.Range("J18").Value = 1.051542595051753
$ws.Range("K18").Value = 1.053363051322834
$ws.Range("L18").Value = 1.046477434453211
$ws.Range("M18").Value = 1.062177384847303

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.045621113021884
$ws.Range("D19").Value = 1.050177746016654
$ws.Range("E19").Value = 1.043280836887256
$ws.Range("F19").Value = 1.059032832850917
$ws.Range("I19").Value = 1.042424707425087
$ws.Range("J19").Value = 1.051634161885123
$ws.Range("K19").Value = 1.053439560178995
$ws.Range("L19").Value = 1.046565732485558
$ws.Range("M19").Value = 1.06226549555885

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.045088972014914
$ws.Range("D20").Value = 1.049771100795816
$ws.Range("E20").Value = 1.042820817898274
$ws.Range("F20").Value = 1.058574772627588
$ws.Range("I20").Value = 1.042278208612338
$ws.Range("J20").Value = 1.051224499533824
$ws.Range("K20").Value = 1.05309722605104
$ws.Range("L20").Value = 1.046170701447976
$ws.Range("M20").Value = 1.061871295906621

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.043356650682585
$ws.Range("D21").Value = 1.048446859104319
$ws.Range("E21").Value = 1.041323468270875
$ws.Range("F21").Value = 1.057083693358386
$ws.Range("I21").Value = 1.041798699925422
$ws.Range("J21").Value = 1.049889916566874
$ws.Range("K21").Value = 1.051981282621639
$ws.Range("L21").Value = 1.044883906572894
$ws.Range("M21").Value = 1.060587093878604

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.04226567677385
$ws.Range("D22").Value = 1.04761253988534
$ws.Range("E22").Value = 1.040380621955986
$ws.Range("F22").Value = 1.056144716671083
$ws.Range("I22").Value = 1.041494725996902
$ws.Range("J22").Value = 1.049048682917633
$ws.Range("K22").Value = 1.05127733024006
$ws.Range("L22").Value = 1.044072889768369
$ws.Range("M22").Value = 1.059777624810362

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.042844195888606
$ws.Range("D23").Value = 1.04805499292112
$ws.Range("E23").Value = 1.040880577979997
$ws.Range("F23").Value = 1.056642628093356
$ws.Range("I23").Value = 1.04165610507359
$ws.Range("J23").Value = 1.049494840808799
$ws.Range("K23").Value = 1.05165072996229
$ws.Range("L23").Value = 1.044503012857238
$ws.Range("M23").Value = 1.060206935112736

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.045117984875067
$ws.Range("D24").Value = 1.049793273237734
$ws.Range("E24").Value = 1.042845897900576
$ws.Range("F24").Value = 1.058599746241455
$ws.Range("I24").Value = 1.042286205673922
$ws.Range("J24").Value = 1.051246838420802
$ws.Range("K24").Value = 1.053115896171094
$ws.Range("L24").Value = 1.046192242027788
$ws.Range("M24").Value = 1.061892791591591

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.047749083448768
$ws.Range("D25").Value = 1.051803176975658
$ws.Range("E25").Value = 1.045120666788838
$ws.Range("F25").Value = 1.060864666435096
$ws.Range("I25").Value = 1.043006662639419
$ws.Range("J25").Value = 1.053270893250965
$ws.Range("K25").Value = 1.054806247793387
$ws.Range("L25").Value = 1.048144189117739
$ws.Range("M25").Value = 1.063840453435576
